# Update the "想去人数" (F column) figures on the "展览" and "全部类型"
# worksheets to the newer snapshot values.

$wb = $excel.ActiveWorkbook

# Map of row -> new value for column F
$updates = @{
    2  = 1576
    3  = 51
    5  = 29
    7  = 2700
    9  = 1721
    10 = 185
    11 = 74
    12 = 586
    13 = 28
    14 = 16
    15 = 99
    16 = 78
    18 = 14
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
